$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The table shrank by one row (old last row 19 removed)
$ws.Range("A19:C19").ClearContents()
